$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 40 : aula 41 - "41. Alerta de sucesso" (simple observation text)
# ---------------------------------------------------------------------------
$ws.Range("B40").Value = 41
$ws.Range("C40").Value = "8. Departamento: Controller & View"
$ws.Range("D40").Value = "41. Alerta de sucesso"
$ws.Range("D40").WrapText = $true

$obsRow40 = "3:32`npara exibir os alerts, é o controller no back-end que envia as variaveis para as páginas, que no caso são mensagens de SUCESS ou FAIL de acordo com o comportamento dos métidos inserir, editar ou excluir. A implementação fica a critério."
$ws.Range("E40").Value = $obsRow40
$ws.Range("E40").WrapText = $true
$ws.Rows.Item(40).RowHeight = 75

# ---------------------------------------------------------------------------
# Row 41 : aula 41 - "41. Alerta de sucesso" (rich-text "IMPORTANTISSIMO" note)
# ---------------------------------------------------------------------------
$ws.Range("B41").Value = 41
$ws.Range("C41").Value = "8. Departamento: Controller & View"
$ws.Range("D41").Value = "41. Alerta de sucesso"
$ws.Range("D41").WrapText = $true

$run1 = "4:50`n"
$run2 = "IMPORTANTISSIMO:"
$run3 = "  ao enviar variaveis EL para o front-end:`nse usar o redirect para retornar para uma pagina, no escopo do método deve usar um parametro ""RedirectAttributes"" passando variavel EL para o front com o metodo ""addFlashAttributes""`nse usar o ModelMap para retornar uma pagina, no escopo do método deve usar um parametro ModelMap passando a variavel EL para o front com o método ""addAttribute"""

$cellE41 = $ws.Range("E41")
$cellE41.Value = $run1 + $run2 + $run3

# base style of the cell: bold, size 12 (this becomes the cell's default font,
# which is what the unformatted first run ("4:50") visually inherits)
$cellE41.Font.Bold = $true
$cellE41.Font.Size = 12
$cellE41.WrapText = $true
$cellE41.Interior.Color = 49407

# "IMPORTANTISSIMO:" in green
$run2Start = $run1.Length + 1
$chars2 = $cellE41.Characters($run2Start, $run2.Length)
$chars2.Font.Bold = $true
$chars2.Font.Size = 12
$chars2.Font.Color = 5287936

# remainder of the note, bold/size 12 (theme text color)
$run3Start = $run2Start + $run2.Length
$chars3 = $cellE41.Characters($run3Start, $run3.Length)
$chars3.Font.Bold = $true
$chars3.Font.Size = 12

$ws.Rows.Item(41).RowHeight = 189

# ---------------------------------------------------------------------------
# Trailing blank rows 42 and 43 (wrap-text style carried down column E)
# ---------------------------------------------------------------------------
$ws.Range("E42").WrapText = $true
$ws.Range("E43").WrapText = $true

# ---------------------------------------------------------------------------
# Update the view so the new rows are visible / selected, as in the diff
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 37
$ws.Range("E41").Select()

Write-Host "done"
